$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F7").Value = 2685
$wsExhibition.Range("F11").Value = 10099
$wsExhibition.Range("F17").Value = 12141

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2685
$wsAll.Range("F12").Value = 10099
$wsAll.Range("F18").Value = 12141
